$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Ang, Bryan -----------------------------------------------------
$ws.Cells.Item(2, 1).Value = "Ang, Bryan"
$ws.Cells.Item(2, 2).Value = "-"
$ws.Cells.Item(2, 3).Value = "-"
$ws.Cells.Item(2, 4).Value = 44115
$ws.Cells.Item(2, 5).Value = "-"
$ws.Cells.Item(2, 6).Value = 3
$ws.Cells.Item(2, 7).Value = "-"
$ws.Cells.Item(2, 8).Value = "zang515@aucklanduni.ac.nz"
$ws.Cells.Item(2, 9).Value = 261940721

# --- Row 3: Wang Ma, Frank --------------------------------------------------
$ws.Cells.Item(3, 1).Value = "Wang Ma, Frank"
$ws.Cells.Item(3, 2).Value = "-"
$ws.Cells.Item(3, 3).Value = "-"
$ws.Cells.Item(3, 4).Value = 44121
$ws.Cells.Item(3, 5).Value = "-"
$ws.Cells.Item(3, 6).Value = 17
$ws.Cells.Item(3, 7).Value = "-"
$ws.Cells.Item(3, 8).Value = "fwan175@aucklanduni.ac.nz"
$ws.Cells.Item(3, 9).Value = 184846458

# --- Rows 4 & 5: clear the old duplicate records, leaving only the ---------
# --- date-formatted placeholder cell in column D ----------------------------
$ws.Range("A4:I4").ClearContents()
$ws.Range("A5:I5").ClearContents()

Write-Output "done"
